$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to Text format so that
# numeric-looking values (e.g. "1.00", "61.628.60") are stored as text,
# matching the inline-string cells used throughout this sheet.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '61.628.60'
$ws.Range("E2").Value = '  +0.30%  '

# Row 3
$ws.Range("D3").Value = '3.443.91'
$ws.Range("E3").Value = '  +2.48%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '577.91'
$ws.Range("E5").Value = '  +0.97%  '

# Row 6
$ws.Range("D6").Value = '147.10'
$ws.Range("E6").Value = '  +7.76%  '

# Row 7
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.442.78'
$ws.Range("E8").Value = '  +2.50%  '

# Row 9
$ws.Range("D9").Value = '0.472'
$ws.Range("E9").Value = '  +0.66%  '

# Row 10
$ws.Range("D10").Value = '7.68'
$ws.Range("E10").Value = '  +2.76%  '

# Row 11
$ws.Range("D11").Value = '0.123'
$ws.Range("E11").Value = '  -0.49%  '

# Row 12
$ws.Range("D12").Value = '0.385'
$ws.Range("E12").Value = '  -0.71%  '

# Row 13
$ws.Range("D13").Value = '4.022.67'
$ws.Range("E13").Value = '  +2.22%  '

# Row 14
$ws.Range("D14").Value = '27.86'
$ws.Range("E14").Value = '  +7.80%  '

# Row 15
$ws.Range("E15").Value = '  -0.62%  '

# Row 16
$ws.Range("D16").Value = '0.0000174'
$ws.Range("E16").Value = '  +0.24%  '

# Row 17
$ws.Range("D17").Value = '3.446.60'
$ws.Range("E17").Value = '  +2.42%  '

# Row 18
$ws.Range("D18").Value = '61.726.98'
$ws.Range("E18").Value = '  +0.32%  '

# Row 19
$ws.Range("D19").Value = '6.25'
$ws.Range("E19").Value = '  +6.83%  '

# Row 20
$ws.Range("D20").Value = '14.05'
$ws.Range("E20").Value = '  +0.60%  '

# Row 21
$ws.Range("D21").Value = '9.35'
$ws.Range("E21").Value = '  +0.49%  '

# Row 22
$ws.Range("D22").Value = '382.69'
$ws.Range("E22").Value = '  +0.98%  '

# Row 23
$ws.Range("D23").Value = '0.564'
$ws.Range("E23").Value = '  +2.18%  '

# Row 24
$ws.Range("D24").Value = '3.592.01'
$ws.Range("E24").Value = '  +2.63%  '

# Row 25
$ws.Range("E25").Value = '  +0.17%  '

# Row 26
$ws.Range("E26").Value = '  +0.44%  '

# Row 27
$ws.Range("D27").Value = '72.17'
$ws.Range("E27").Value = '  +1.32%  '

# Row 28
$ws.Range("E28").Value = '  -1.30%  '

# Row 29
$ws.Range("D29").Value = '0.178'
$ws.Range("E29").Value = '  +8.38%  '

# Row 30
$ws.Range("D30").Value = '7.70'
$ws.Range("E30").Value = '  +2.85%  '

# Row 31
$ws.Range("E31").Value = '  -12.06%  '

# Row 32
$ws.Range("E32").Value = '  -0.04%  '

# Row 33
$ws.Range("D33").Value = '8.18'
$ws.Range("E33").Value = '  +0.08%  '

# Row 34
$ws.Range("E34").Value = '  +1.07%  '

# Row 35
$ws.Range("E35").Value = '  -0.05%  '

# Row 36
$ws.Range("D36").Value = '23.98'
$ws.Range("E36").Value = '  +1.78%  '

# Row 37
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = '7.01'
$ws.Range("E37").Value = '  +3.26%  '

# Row 38
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '5.20'
$ws.Range("E38").Value = '  -0.07%  '

# Row 39
$ws.Range("D39").Value = '1.55'
$ws.Range("E39").Value = '  +2.13%  '

# Row 40
$ws.Range("D40").Value = '165.54'
$ws.Range("E40").Value = '  -0.01%  '

# Row 41
$ws.Range("D41").Value = '0.0781'
$ws.Range("E41").Value = '  +2.56%  '

# Row 42
$ws.Range("D42").Value = '25.82'
$ws.Range("E42").Value = '  +8.01%  '

# Row 43
$ws.Range("D43").Value = '0.791'
$ws.Range("E43").Value = '  +2.74%  '

# Row 44
$ws.Range("E44").Value = '  +0.02%  '

# Row 45
$ws.Range("D45").Value = '1.72'
$ws.Range("E45").Value = '  +0.30%  '

# Row 46
$ws.Range("D46").Value = '4.46'
$ws.Range("E46").Value = '  +1.84%  '

# Row 47
$ws.Range("D47").Value = '42.04'
$ws.Range("E47").Value = '  +1.51%  '

# Row 48
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.617.00'
$ws.Range("E48").Value = '  +10.61%  '

# Row 49
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value = '1.17'
$ws.Range("E49").Value = '  -2.73%  '

# Row 50
$ws.Range("D50").Value = '23.54'
$ws.Range("E50").Value = '  +2.32%  '

# Row 51
$ws.Range("D51").Value = '6.86'
$ws.Range("E51").Value = '  +0.33%  '
